$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = '"COVID-19\u2013Associated Hepatitis Complicating Recent Living Donor Liver Transplantation"'
$ws.Range("D2").Value = '"We present a case of COVID-19 hepatitis in a living donor liver allograft recipient whose donor subsequently tested positive for COVID-19. The patient is a female infant with biliary atresia (failed Kasai procedure).
 She recovered well, with improving liver function tests for 4 days.
 On postoperative day 4 the patient developed respiratory distress and fever.
 COVID-19 testing (polymerase chain reaction) was positive.
 Liver function test results increased approximately 5-fold.
 Liver biopsy showed moderate acute hepatitis with prominent clusters of apoptotic hepatocytes and associated cellular debris.
 Lobular lymphohistiocytic inflammation was noted.
 Typical portal features of mild to moderate acute cellular rejection were also noted.
'
$ws.Range("E2").Value = '[Stephen M.%Lagana%xref no email%1,  Simona%De Michele%xref no email%1,  Michael J.%Lee%xref no email%1,  Jean C.%Emond%xref no email%1,  Adam D.%Griesemer%xref no email%1,  Sheryl A.%Tulin-Silver%xref no email%1,  Elizabeth C.%Verna%xref no email%1,  Mercedes%Martinez%xref no email%1,  Jay H.%Lefkowitch%xref no email%1]'
$ws.Range("F2").Value = '10.5858/arpa.2020-0186-SA'
$ws.Range("G2").Value = 'CROSSREF'
$ws.Range("H2").Value = '2023-05-11'
$ws.Range("E3").Value = '[Mario%Fernández-Ruiz%NULL%2,         Amado%Andrés%NULL%2,         Carmelo%Loinaz%NULL%2,         Juan F.%Delgado%NULL%2,         Francisco%López-Medrano%NULL%2,         Rafael%San Juan%NULL%2,         Esther%González%NULL%2,         Natalia%Polanco%NULL%2,         María D.%Folgueira%NULL%2,         Antonio%Lalueza%NULL%2,         Carlos%Lumbreras%NULL%2,         José M.%Aguado%NULL%2]'
$ws.Range("E4").Value = '[David%Marx%NULL%1,         Bruno%Moulin%NULL%1,         Samira%Fafi-Kremer%NULL%2,         Ilies%Benotmane%NULL%1,         Gabriela%Gautier%NULL%1,         Peggy%Perrin%NULL%1,         Sophie%Caillard%NULL%1]'
$ws.Range("E5").Value = '[Marilù%Bartiromo%bartiromom@aou-careggi.toscana.it%1,         Beatrice%Borchi%NULL%1,         Annarita%Botta%NULL%2,         Annarita%Botta%NULL%0,         Alfredo%Bagalà%NULL%2,         Alfredo%Bagalà%NULL%0,         Gianmarco%Lugli%NULL%1,         Marta%Tilli%NULL%2,         Marta%Tilli%NULL%0,         Annalisa%Cavallo%NULL%1,         Brunilda%Xhaferi%NULL%1,         Roberta%Cutruzzulà%NULL%1,         Augusto%Vaglio%NULL%1,         Silvia%Bresci%NULL%1,         Aida%Larti%NULL%1,         Alessandro%Bartoloni%NULL%1,         Calogero%Cirami%NULL%1]'
$ws.Range("E6").Value = '[Ling%Ning%NULL%1,         Lei%Liu%NULL%0,         Wenyuan%Li%NULL%1,         Hongtao%Liu%NULL%1,         Jizhou%Wang%NULL%1,         Ziqin%Yao%NULL%1,         Shengyu%Zhang%NULL%1,         Desheng%Zhao%NULL%1,         Björn%Nashan%NULL%1,         Aizong%Shen%NULL%1,         Lianxin%Liu%NULL%2,         Lei%Li%NULL%3]'
$ws.Range("E7").Value = '[Junpeng%Wang%NULL%1,         Xin%Li%NULL%3,         Guanghui%Cao%NULL%1,         Xiaoqiang%Wu%NULL%1,         Zhiwei%Wang%NULL%1,         Tianzhong%Yan%ytz460@hotmail.com%1]'
$ws.Range("E8").Value = '[Hui%Zhang%NULL%1,         Yan%Chen%NULL%1,         Quan%Yuan%NULL%1,         Qiu-Xiang%Xia%NULL%1,         Xian-Peng%Zeng%NULL%1,         Jing-Tao%Peng%NULL%1,         Jing%Liu%NULL%3,         Xing-Yuan%Xiao%NULL%1,         Guo-Song%Jiang%NULL%1,         Han-Yu%Xiao%NULL%1,         Liang-Bo%Xie%NULL%1,         Jing%Chen%NULL%2,         Jia-Li%Liu%NULL%1,         Xiong%Xiao%NULL%1,         Hua%Su%NULL%1,         Chun%Zhang%NULL%1,         Xiao-Ping%Zhang%NULL%1,         Hua%Yang%NULL%1,         Heng%Li%NULL%1,         Zhen-Di%Wang%wangzhendi@gmail.com%1]'
$ws.Range("E9").Value = '[Song%Chen%NULL%1,         Qin%Yin%NULL%1,         Huibo%Shi%NULL%2,         Dunfeng%Du%NULL%1,         Sheng%Chang%NULL%1,         Li%Ni%NULL%1,         Haifang%Qiu%NULL%1,         Zhishui%Chen%NULL%2,         Jixian%Zhang%NULL%1,         Weijie%Zhang%NULL%1]'
$ws.Range("E10").Value = '[Elena%Seminari%NULL%1,         Marta%Colaneri%NULL%1,         Margherita%Sambo%NULL%1,         Ilaria%Gallazzi%NULL%1,         Angela%Di Matteo%NULL%1,         Silvia%Roda%NULL%1,         Raffaele%Bruno%NULL%0,         NULL%NULL%NULL%0]'
$ws.Range("E11").Value = '[Ilaria%Gandolfini%NULL%1,         Marco%Delsante%NULL%1,         Enrico%Fiaccadori%NULL%1,         Gianluigi%Zaza%NULL%1,         Lucio%Manenti%NULL%1,         Anna%Degli Antoni%NULL%1,         Licia%Peruzzi%NULL%1,         Leonardo V.%Riella%NULL%1,         Paolo%Cravedi%NULL%1,         Umberto%Maggiore%NULL%1]'
$ws.Range("E12").Value = '[Elena%Guillen%NULL%1,         Gaston J.%Pineiro%NULL%1,         Ignacio%Revuelta%NULL%1,         Diana%Rodriguez%NULL%1,         Marta%Bodro%NULL%1,         Asunción%Moreno%NULL%1,         Josep M.%Campistol%NULL%1,         Fritz%Diekmann%NULL%1,         Pedro%Ventura-Aguiar%NULL%1]'
$ws.Range("E13").Value = '[Mario%Fernández-Ruiz%NULL%0,         Amado%Andrés%NULL%0,         Carmelo%Loinaz%NULL%0,         Juan F.%Delgado%NULL%0,         Francisco%López-Medrano%NULL%0,         Rafael%San Juan%NULL%0,         Esther%González%NULL%0,         Natalia%Polanco%NULL%0,         María D.%Folgueira%NULL%0,         Antonio%Lalueza%NULL%0,         Carlos%Lumbreras%NULL%0,         José M.%Aguado%NULL%0]'
$ws.Range("E14").Value = '[Zibiao%Zhong%NULL%1,         Qiuyan%Zhang%NULL%1,         Haoyang%Xia%NULL%1,         Aiping%Wang%NULL%1,         Wenjin%Liang%NULL%1,         Wei%Zhou%NULL%1,         Lihua%Zhou%NULL%1,         Xiao%Liu%NULL%1,         Lingzhang%Rao%NULL%1,         Zhifeng%Li%NULL%1,         Zhiyong%Peng%NULL%1,         Pingzheng%Mo%NULL%0,         Yong%Xiong%NULL%0,         Shaojun%Ye%NULL%1,         Yanfeng%Wang%NULL%1,         Qifa%Ye%NULL%1]'
$ws.Range("C15").Value = '"Perioperative Presentation of COVID\u201019 Disease in a Liver Transplant Recipient"'
$ws.Range("E15").Value = '[Juanjuan%Qin%xref no email%1,  Haitao%Wang%xref no email%1,  Xian%Qin%xref no email%1,  Peng%Zhang%xref no email%1,  Lihua%Zhu%xref no email%1,  Jingjing%Cai%xref no email%1,  Yufeng%Yuan%xref no email%1,  Hongliang%Li%xref no email%1]'
$ws.Range("F15").Value = '10.1002/hep.31257'
$ws.Range("H15").Value = '2023-04-24'
$ws.Range("E16").Value = '[Jiao-Feng%Huang%NULL%1,         Kenneth I.%Zheng%NULL%1,         Jacob%George%NULL%1,         Hai-Nv%Gao%NULL%1,         Ru-Nan%Wei%NULL%1,         Hua-Dong%Yan%NULL%1,         Ming-Hua%Zheng%NULL%2]'
$ws.Range("E17").Value = '[Bin%Liu%NULL%1,         Yangzhong%Wang%NULL%1,         Yuanyuan%Zhao%NULL%1,         Huibo%Shi%NULL%0,         Fanjun%Zeng%NULL%1,         Zhishui%Chen%NULL%0]'
$ws.Range("E18").Value = '[Sherrie%Bhoori%NULL%1,         Roberta Elisa%Rossi%NULL%1,         Davide%Citterio%NULL%1,         Vincenzo%Mazzaferro%NULL%1]'
